$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.09389020502567291
$ws.Cells.Item(2, 2).Value = 0.9777096509933472
$ws.Cells.Item(2, 3).Value = 0.02287952229380608
$ws.Cells.Item(2, 4).Value = 0.9943627715110779

$ws.Cells.Item(3, 1).Value = 0.02444709837436676
$ws.Cells.Item(3, 2).Value = 0.9978806972503662
$ws.Cells.Item(3, 3).Value = 0.01488614827394485
$ws.Cells.Item(3, 4).Value = 0.998120903968811

$ws.Cells.Item(4, 1).Value = 0.01320934575051069
$ws.Cells.Item(4, 2).Value = 0.9981456398963928
$ws.Cells.Item(4, 3).Value = 0.01363444328308105
$ws.Cells.Item(4, 4).Value = 0.9984477162361145

$ws.Cells.Item(5, 1).Value = 0.008361239917576313
$ws.Cells.Item(5, 2).Value = 0.9983726739883423
$ws.Cells.Item(5, 3).Value = 0.07270972430706024
$ws.Cells.Item(5, 4).Value = 0.9754084944725037

$ws.Cells.Item(6, 1).Value = 0.006124851759523153
$ws.Cells.Item(6, 2).Value = 0.9985997676849365
$ws.Cells.Item(6, 3).Value = 0.002859361469745636
$ws.Cells.Item(6, 4).Value = 0.998774528503418

$ws.Cells.Item(7, 1).Value = 0.006062730215489864
$ws.Cells.Item(7, 2).Value = 0.998505175113678
$ws.Cells.Item(7, 3).Value = 0.002121088793501258
$ws.Cells.Item(7, 4).Value = 0.9989379048347473

$ws.Cells.Item(8, 1).Value = 0.006021596491336823
$ws.Cells.Item(8, 2).Value = 0.9984673261642456
$ws.Cells.Item(8, 3).Value = 0.002340758219361305
$ws.Cells.Item(8, 4).Value = 0.998774528503418

$ws.Cells.Item(9, 1).Value = 0.005663695745170116
$ws.Cells.Item(9, 2).Value = 0.9985240697860718
$ws.Cells.Item(9, 3).Value = 0.00321565568447113
$ws.Cells.Item(9, 4).Value = 0.9989379048347473

$ws.Cells.Item(10, 1).Value = 0.004984347615391016
$ws.Cells.Item(10, 2).Value = 0.9985619187355042
$ws.Cells.Item(10, 3).Value = 0.002122032456099987
$ws.Cells.Item(10, 4).Value = 0.9989379048347473

$ws.Cells.Item(11, 1).Value = 0.004646725486963987
$ws.Cells.Item(11, 2).Value = 0.9986565113067627
$ws.Cells.Item(11, 3).Value = 0.003982760943472385
$ws.Cells.Item(11, 4).Value = 0.9989379048347473

$ws.Cells.Item(12, 1).Value = 0.004548538941890001
$ws.Cells.Item(12, 2).Value = 0.9986186623573303
$ws.Cells.Item(12, 3).Value = 0.001969551201909781
$ws.Cells.Item(12, 4).Value = 0.9989379048347473

$ws.Cells.Item(13, 1).Value = 0.004830501042306423
$ws.Cells.Item(13, 2).Value = 0.9985429644584656
$ws.Cells.Item(13, 3).Value = 0.002122444799169898
$ws.Cells.Item(13, 4).Value = 0.9989379048347473

$ws.Cells.Item(14, 1).Value = 0.005139888729900122
$ws.Cells.Item(14, 2).Value = 0.998505175113678
$ws.Cells.Item(14, 3).Value = 0.02295811474323273
$ws.Cells.Item(14, 4).Value = 0.9929738640785217

$ws.Cells.Item(15, 1).Value = 0.004468801431357861
$ws.Cells.Item(15, 2).Value = 0.9986943602561951
$ws.Cells.Item(15, 3).Value = 0.001441582222469151
$ws.Cells.Item(15, 4).Value = 0.9989379048347473

$ws.Cells.Item(16, 1).Value = 0.004495252389460802
$ws.Cells.Item(16, 2).Value = 0.998751163482666
$ws.Cells.Item(16, 3).Value = 0.001996832434087992
$ws.Cells.Item(16, 4).Value = 0.9990196228027344

$ws.Cells.Item(17, 1).Value = 0.004978906363248825
$ws.Cells.Item(17, 2).Value = 0.9986565113067627
$ws.Cells.Item(17, 3).Value = 0.001904443488456309
$ws.Cells.Item(17, 4).Value = 0.9994280934333801

$ws.Cells.Item(18, 1).Value = 0.005232904106378555
$ws.Cells.Item(18, 2).Value = 0.9987700581550598
$ws.Cells.Item(18, 3).Value = 0.001857162569649518
$ws.Cells.Item(18, 4).Value = 0.9990196228027344

$ws.Cells.Item(19, 1).Value = 0.004224963020533323
$ws.Cells.Item(19, 2).Value = 0.9989971518516541
$ws.Cells.Item(19, 3).Value = 0.00204593944363296
$ws.Cells.Item(19, 4).Value = 0.9996731877326965

$ws.Cells.Item(20, 1).Value = 0.004253908526152372
$ws.Cells.Item(20, 2).Value = 0.9989781975746155
$ws.Cells.Item(20, 3).Value = 0.00197804975323379
$ws.Cells.Item(20, 4).Value = 0.9996731877326965

$ws.Cells.Item(21, 1).Value = 0.004611051641404629
$ws.Cells.Item(21, 2).Value = 0.9990160465240479
$ws.Cells.Item(21, 3).Value = 0.001526396372355521
$ws.Cells.Item(21, 4).Value = 0.9996731877326965

$ws.Cells.Item(22, 1).Value = 0.004608353599905968
$ws.Cells.Item(22, 2).Value = 0.9989971518516541
$ws.Cells.Item(22, 3).Value = 0.002244416391476989
$ws.Cells.Item(22, 4).Value = 0.9996731877326965

$ws.Cells.Item(23, 1).Value = 0.004847287200391293
$ws.Cells.Item(23, 2).Value = 0.9989971518516541
$ws.Cells.Item(23, 3).Value = 0.002315493067726493
$ws.Cells.Item(23, 4).Value = 0.9995098114013672

$ws.Cells.Item(24, 1).Value = 0.005048518534749746
$ws.Cells.Item(24, 2).Value = 0.9987322092056274
$ws.Cells.Item(24, 3).Value = 0.002286312403157353
$ws.Cells.Item(24, 4).Value = 0.9996731877326965

$ws.Cells.Item(25, 1).Value = 0.00383491488173604
$ws.Cells.Item(25, 2).Value = 0.9991484880447388
$ws.Cells.Item(25, 3).Value = 0.002175794914364815
$ws.Cells.Item(25, 4).Value = 0.9996731877326965

$ws.Cells.Item(26, 1).Value = 0.005678361281752586
$ws.Cells.Item(26, 2).Value = 0.9988646507263184
$ws.Cells.Item(26, 3).Value = 0.00106322206556797
$ws.Cells.Item(26, 4).Value = 0.9998366236686707

$ws.Cells.Item(27, 1).Value = 0.004518436267971992
$ws.Cells.Item(27, 2).Value = 0.9990538954734802
$ws.Cells.Item(27, 3).Value = 0.002096576150506735
$ws.Cells.Item(27, 4).Value = 0.9996731877326965

$ws.Cells.Item(28, 1).Value = 0.003697582520544529
$ws.Cells.Item(28, 2).Value = 0.9992809295654297
$ws.Cells.Item(28, 3).Value = 0.001654272782616317
$ws.Cells.Item(28, 4).Value = 0.9996731877326965

$ws.Cells.Item(29, 1).Value = 0.004185241181403399
$ws.Cells.Item(29, 2).Value = 0.9991863369941711
$ws.Cells.Item(29, 3).Value = 0.001593250082805753
$ws.Cells.Item(29, 4).Value = 0.9998366236686707

$ws.Cells.Item(30, 1).Value = 0.003962030634284019
$ws.Cells.Item(30, 2).Value = 0.9992809295654297
$ws.Cells.Item(30, 3).Value = 0.001081951893866062
$ws.Cells.Item(30, 4).Value = 0.9995915293693542

$ws.Cells.Item(31, 1).Value = 0.00411638617515564
$ws.Cells.Item(31, 2).Value = 0.9992998838424683
$ws.Cells.Item(31, 3).Value = 0.001442144624888897
$ws.Cells.Item(31, 4).Value = 0.9995915293693542

$ws.Cells.Item(32, 1).Value = 0.004763389937579632
$ws.Cells.Item(32, 2).Value = 0.9992431402206421
$ws.Cells.Item(32, 3).Value = 0.001249335240572691
$ws.Cells.Item(32, 4).Value = 0.9996731877326965

$ws.Cells.Item(33, 1).Value = 0.004428584594279528
$ws.Cells.Item(33, 2).Value = 0.9992241859436035
$ws.Cells.Item(33, 3).Value = 0.00166610348969698
$ws.Cells.Item(33, 4).Value = 0.9996731877326965

$ws.Cells.Item(34, 1).Value = 0.00376812880858779
$ws.Cells.Item(34, 2).Value = 0.9993377327919006
$ws.Cells.Item(34, 3).Value = 0.0007695319363847375
$ws.Cells.Item(34, 4).Value = 0.9997549057006836

$ws.Cells.Item(35, 1).Value = 0.004741899203509092
$ws.Cells.Item(35, 2).Value = 0.9991484880447388
$ws.Cells.Item(35, 3).Value = 0.001457889680750668
$ws.Cells.Item(35, 4).Value = 0.9996731877326965

$ws.Cells.Item(36, 1).Value = 0.004390857648104429
$ws.Cells.Item(36, 2).Value = 0.9992620348930359
$ws.Cells.Item(36, 3).Value = 0.001254474860616028
$ws.Cells.Item(36, 4).Value = 0.9997549057006836

$ws.Cells.Item(37, 1).Value = 0.004552315920591354
$ws.Cells.Item(37, 2).Value = 0.9992431402206421
$ws.Cells.Item(37, 3).Value = 0.002543745562434196
$ws.Cells.Item(37, 4).Value = 0.9995915293693542

$ws.Cells.Item(38, 1).Value = 0.00490929139778018
$ws.Cells.Item(38, 2).Value = 0.999129593372345
$ws.Cells.Item(38, 3).Value = 0.002570071490481496
$ws.Cells.Item(38, 4).Value = 0.9995915293693542

$ws.Cells.Item(39, 1).Value = 0.004172834567725658
$ws.Cells.Item(39, 2).Value = 0.9992809295654297
$ws.Cells.Item(39, 3).Value = 0.002720575081184506
$ws.Cells.Item(39, 4).Value = 0.9996731877326965

$ws.Cells.Item(40, 1).Value = 0.00468365428969264
$ws.Cells.Item(40, 2).Value = 0.9992052912712097
$ws.Cells.Item(40, 3).Value = 0.002715761307626963
$ws.Cells.Item(40, 4).Value = 0.9996731877326965

$ws.Cells.Item(41, 1).Value = 0.005102598108351231
$ws.Cells.Item(41, 2).Value = 0.9992052912712097
$ws.Cells.Item(41, 3).Value = 0.0007813895354047418
$ws.Cells.Item(41, 4).Value = 0.9999182820320129

$ws.Cells.Item(42, 1).Value = 0.004313069861382246
$ws.Cells.Item(42, 2).Value = 0.9992620348930359
$ws.Cells.Item(42, 3).Value = 0.0006634125020354986
$ws.Cells.Item(42, 4).Value = 1

$ws.Cells.Item(43, 1).Value = 0.004116281401365995
$ws.Cells.Item(43, 2).Value = 0.9993187785148621
$ws.Cells.Item(43, 3).Value = 0.1140810921788216
$ws.Cells.Item(43, 4).Value = 0.9838235378265381

$ws.Cells.Item(44, 1).Value = 0.004473397973924875
$ws.Cells.Item(44, 2).Value = 0.9992241859436035
$ws.Cells.Item(44, 3).Value = 0.0007821788312867284
$ws.Cells.Item(44, 4).Value = 1

$ws.Cells.Item(45, 1).Value = 0.004450363572686911
$ws.Cells.Item(45, 2).Value = 0.9992241859436035
$ws.Cells.Item(45, 3).Value = 0.00222443975508213
$ws.Cells.Item(45, 4).Value = 0.9996731877326965

$ws.Cells.Item(46, 1).Value = 0.00369930500164628
$ws.Cells.Item(46, 2).Value = 0.999375581741333
$ws.Cells.Item(46, 3).Value = 0.002451557666063309
$ws.Cells.Item(46, 4).Value = 0.9996731877326965

$ws.Cells.Item(47, 1).Value = 0.003945261240005493
$ws.Cells.Item(47, 2).Value = 0.9993566274642944
$ws.Cells.Item(47, 3).Value = 0.002171063795685768
$ws.Cells.Item(47, 4).Value = 0.9996731877326965

$ws.Cells.Item(48, 1).Value = 0.006097278092056513
$ws.Cells.Item(48, 2).Value = 0.9990349411964417
$ws.Cells.Item(48, 3).Value = 0.001885453704744577
$ws.Cells.Item(48, 4).Value = 0.9996731877326965

$ws.Cells.Item(49, 1).Value = 0.003999381326138973
$ws.Cells.Item(49, 2).Value = 0.9993377327919006
$ws.Cells.Item(49, 3).Value = 0.001836212119087577
$ws.Cells.Item(49, 4).Value = 0.9996731877326965

$ws.Cells.Item(50, 1).Value = 0.004657507874071598
$ws.Cells.Item(50, 2).Value = 0.9990160465240479
$ws.Cells.Item(50, 3).Value = 0.001600455143488944
$ws.Cells.Item(50, 4).Value = 0.9996731877326965

$ws.Cells.Item(51, 1).Value = 0.004102061036974192
$ws.Cells.Item(51, 2).Value = 0.999072790145874
$ws.Cells.Item(51, 3).Value = 0.001546047045849264
$ws.Cells.Item(51, 4).Value = 0.9997549057006836
